$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing content first - the table layout is being replaced.
$ws.UsedRange.ClearContents()

# Write cells in the same order the original authoring session introduced new
# shared strings, so the rebuilt sharedStrings table lines up with the source
# edit (section headers, then the new "Within" column, then the crew table).
$ws.Range("A1").Value = "** Crew"
$ws.Range("A6").Value = "** Habitat"
$ws.Range("A12").Value = "Potato"
$ws.Range("F2").Value = "Within"
$ws.Range("B2").Value = "Age"
$ws.Range("C2").Value = "Weight"
$ws.Range("D2").Value = "Height"
$ws.Range("E2").Value = "Sex"
$ws.Range("A3").Value = "astro1"
$ws.Range("E3").Value = "M"
$ws.Range("A4").Value = "astro2"
$ws.Range("E4").Value = "F"

# ---- remaining "** Crew" table cells (rows 1-4) ----
$ws.Range("A2").Value = "Name"
$ws.Range("B3").Value = 35
$ws.Range("C3").Value = 75
$ws.Range("D3").Value = 1.7649999999999999
$ws.Range("F3").Value = "Hab"
$ws.Range("B4").Value = 35
$ws.Range("C4").Value = 55
$ws.Range("D4").Value = 1.63
$ws.Range("F4").Value = "Hab"

# ---- remaining "** Habitat" table cells (rows 6-8) ----
$ws.Range("A7").Value = "Name"
$ws.Range("B7").Value = "Temperature"
$ws.Range("C7").Value = "Volume"
$ws.Range("D7").Value = "O2Percentage"
$ws.Range("E7").Value = "CO2Percentage"
$ws.Range("F7").Value = "N2Percentage"
$ws.Range("G7").Value = "H2OPercentage"
$ws.Range("H7").Value = "leakPercentage"
$ws.Range("I7").Value = "Contains"

$ws.Range("A8").Value = "Hab"
$ws.Range("B8").Value = 70.3
$ws.Range("C8").Value = 25000
$ws.Range("D8").Value = 0.26500000000000001
$ws.Range("E8").Value = 0.001
$ws.Range("F8").Value = 0.73399999999999999
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = "0.05/24"
$ws.Range("I8").Value = "LettuceReactor"

# ---- remaining "** Crop" table cells (rows 10-12) ----
$ws.Range("A10").Value = "** Crop"
$ws.Range("A11").Value = "Name"
$ws.Range("B11").Value = "Params"
$ws.Range("C11").Value = "Within"

$ws.Range("B12").Value = "/Users/desho/echusOverlook/Simulation/MarsOne/bean_test.xlsx"
$ws.Range("C12").Value = "Hab"

# Match the saved selection/active cell from the authored edit.
$ws.Range("F18").Select()
